# Insert a new "chemical_recycling_pyrolysis" parameter row right after the
# existing "chemical_recycling_gasification" row (currently row 9), pushing
# every row from the old row 10 ("fossil_routes") onward down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new, blank row at position 10 (shifts rows 10-24 down to 11-25).
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row with the new parameter/value pair.
$ws.Range("A10").Value = "chemical_recycling_pyrolysis"
$ws.Range("B10").Value = $true

Write-Output "Inserted chemical_recycling_pyrolysis row at row 10"
